{"js": "// Update the date line and the 20x5 table of arithmetic answers.\n// Values are applied strictly by document order (row-major, top-to-bottom,\n// left-to-right) so that a repeated \"before\" string (e.g. \"55+24=79\" occurs\n// twice) still maps to the correct distinct replacement.\n\nconst newValues = [\n  \"27+24=51\", \"27+14=41\", \"54+15=69\", \"68-66=2\", \"75-37=38\",\n  \"95-55=40\", \"21+62=83\", \"4+79=83\", \"41+13=54\", \"71-34=37\",\n  \"55-0=55\", \"32+61=93\", \"83-32=51\", \"52+8=60\", \"73-37=36\",\n  \"50+23=73\", \"65+20=85\", \"33+33=66\", \"33+42=75\", \"36+41=77\",\n  \"79-57=22\", \"87-4=83\", \"21+4=25\", \"79-40=39\", \"18-7=11\",\n  \"39-10=29\", \"68-18=50\", \"54+45=99\", \"56+25=81\", \"38+60=98\",\n  \"52+8=60\", \"53+36=89\", \"77+10=87\", \"88-13=75\", \"4+7=11\",\n  \"31+3=34\", \"92-72=20\", \"24+17=41\", \"59+31=90\", \"4+28=32\",\n  \"9+8=17\", \"0+45=45\", \"0+10=10\", \"37+58=95\", \"91-51=40\",\n  \"72-24=48\", \"89-17=72\", \"89-34=55\", \"28+2=30\", \"19+76=95\",\n  \"77-51=26\", \"47+11=58\", \"56-18=38\", \"51+43=94\", \"84-67=17\",\n  \"84-26=58\", \"7+63=70\", \"70-39=31\", \"96-93=3\", \"26+62=88\",\n  \"95-61=34\", \"97-40=57\", \"52+37=89\", \"62+35=97\", \"42+29=71\",\n  \"6+59=65\", \"80+14=94\", \"72-20=52\", \"19+55=74\", \"48+46=94\",\n  \"75-68=7\", \"25-3=22\", \"28+20=48\", \"27+7=34\", \"52-41=11\",\n  \"17+5=22\", \"63-30=33\", \"29+0=29\", \"89-63=26\", \"28+54=82\",\n  \"41+33=74\", \"53+13=66\", \"90-84=6\", \"77-9=68\", \"57-26=31\",\n  \"0+86=86\", \"24+72=96\", \"32+36=68\", \"24+14=38\", \"35+54=89\",\n  \"7+56=63\", \"5+59=64\", \"23+50=73\", \"37-34=3\", \"92-75=17\",\n  \"78-42=36\", \"85-79=6\", \"60+26=86\", \"39+47=86\", \"40+16=56\",\n];\n\n// 1) Update the date paragraph (first paragraph of the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\nif (dateParagraph.text.trim() === \"2023-05-19 Friday\") {\n  dateParagraph.insertText(\"2023-05-20 Saturday\", \"Replace\");\n}\n\n// 2) Update every cell of the (single) answers table, in row-major order.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = 5; // fixed 5-column grid per the document's tblGrid\nlet k = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    if (k >= newValues.length) break;\n    const cell = table.getCell(r, c);\n    cell.value = newValues[k];\n    k++;\n  }\n}\nawait context.sync();\n", "ps1": "# Update the date line and the 20x5 table of arithmetic answers.\n# Values are applied strictly by document order (row-major, top-to-bottom,\n# left-to-right) so that a repeated \"before\" string (e.g. \"55+24=79\" occurs\n# twice) still maps to the correct distinct replacement.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph (first paragraph of the body).\n$p1 = $d.Paragraphs(1)\nif ($p1.Range.Text.TrimEnd(\"`r\") -eq \"2023-05-19 Friday\") {\n    $p1.Range.Text = \"2023-05-20 Saturday\"\n}\n\n# 2) Update every cell of the (single) answers table, in row-major order.\n$newValues = @(\n    \"27+24=51\",\"27+14=41\",\"54+15=69\",\"68-66=2\",\"75-37=38\",\n    \"95-55=40\",\"21+62=83\",\"4+79=83\",\"41+13=54\",\"71-34=37\",\n    \"55-0=55\",\"32+61=93\",\"83-32=51\",\"52+8=60\",\"73-37=36\",\n    \"50+23=73\",\"65+20=85\",\"33+33=66\",\"33+42=75\",\"36+41=77\",\n    \"79-57=22\",\"87-4=83\",\"21+4=25\",\"79-40=39\",\"18-7=11\",\n    \"39-10=29\",\"68-18=50\",\"54+45=99\",\"56+25=81\",\"38+60=98\",\n    \"52+8=60\",\"53+36=89\",\"77+10=87\",\"88-13=75\",\"4+7=11\",\n    \"31+3=34\",\"92-72=20\",\"24+17=41\",\"59+31=90\",\"4+28=32\",\n    \"9+8=17\",\"0+45=45\",\"0+10=10\",\"37+58=95\",\"91-51=40\",\n    \"72-24=48\",\"89-17=72\",\"89-34=55\",\"28+2=30\",\"19+76=95\",\n    \"77-51=26\",\"47+11=58\",\"56-18=38\",\"51+43=94\",\"84-67=17\",\n    \"84-26=58\",\"7+63=70\",\"70-39=31\",\"96-93=3\",\"26+62=88\",\n    \"95-61=34\",\"97-40=57\",\"52+37=89\",\"62+35=97\",\"42+29=71\",\n    \"6+59=65\",\"80+14=94\",\"72-20=52\",\"19+55=74\",\"48+46=94\",\n    \"75-68=7\",\"25-3=22\",\"28+20=48\",\"27+7=34\",\"52-41=11\",\n    \"17+5=22\",\"63-30=33\",\"29+0=29\",\"89-63=26\",\"28+54=82\",\n    \"41+33=74\",\"53+13=66\",\"90-84=6\",\"77-9=68\",\"57-26=31\",\n    \"0+86=86\",\"24+72=96\",\"32+36=68\",\"24+14=38\",\"35+54=89\",\n    \"7+56=63\",\"5+59=64\",\"23+50=73\",\"37-34=3\",\"92-75=17\",\n    \"78-42=36\",\"85-79=6\",\"60+26=86\",\"39+47=86\",\"40+16=56\"\n)\n\n$t = $d.Tables(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n$k = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($k -ge $newValues.Count) { continue }\n        $t.Cell($r, $c).Range.Text = $newValues[$k]\n        $k++\n    }\n}\n"}
